$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Quantidade"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Projeto"
$ws.Range("E1").Value = "Técnico"
$ws.Range("F1").Value = "IDGEO"

# Data rows
$data = @(
    @("HASTE INOX 1M", 25, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("PONTEIRA FIXA INOX", 6, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("BOMBA PNEUMATICA", 4, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("TE DE INJEÇÃO", 4, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("MISTURADOR", 3, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("T INOX COM SAIDA MANOMETRO", 4, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("SACADOR DE HASTE", 1, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045"),
    @("BATEDOR", 1, "uso", "QUALITA", "AMAURI / THIAGO SILVA", "PR24045")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
